# Insert a new weekly price record for "Vega Monumental Concepción" (Mango)
# at row 107, shifting the existing rows 107:179 down to 108:180.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 107 (pushes old 107..179 -> 108..180)
$ws.Rows.Item(107).Insert()

# Populate the new row 107 with the latest weekly observation
$ws.Range("A107").Value = 11
$ws.Range("B107").Value = "Vega Monumental Concepción"
$ws.Range("C107").Value = "Bíobío"
$ws.Range("D107").Value = 45126
$ws.Range("E107").Value = 8
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100108
$ws.Range("H107").Value = "Tropicales y subtropicales"
$ws.Range("I107").Value = 100108002
$ws.Range("J107").Value = "Mango"
$ws.Range("K107").Value = "Sin especificar"
$ws.Range("L107").Value = "Primera"
$ws.Range("M107").Value = 100
$ws.Range("N107").Value = 7500
$ws.Range("O107").Value = 8000
$ws.Range("P107").Value = 7750
$ws.Range("Q107").Value = "$/bandeja 4 kilos"
$ws.Range("R107").Value = "Brasil"
$ws.Range("S107").Value = 1938
$ws.Range("T107").Value = 4
